$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 8 (05-12-2015): I8 gets 0.5 (overtime hours) ---
$ws.Range("I8").Value = 0.5

# --- 2. Row 9 (05-13-2015): turn from an "absent" (red) day into an
#        "official business" (blue) day: copy the blue fill from row 8,
#        then populate the OB time columns + remarks ---
$ws.Range("A9:P9").Interior.Color = $ws.Range("A8:P8").Interior.Color

$ws.Range("K9").Value = "08:00:00"
$ws.Range("L9").Value = "08:30:00"
$ws.Range("M9").Value = "18:30:00"
$ws.Range("N9").Value = "18:30:00"
$ws.Range("P9").Value = "~OB Others|PDIS Support due to software corruption.| R"

# --- 3. Remarks text tweaks: append " R " marker / " ~ = " infix ---
$ws.Range("P5").Value = "`" ~OT ~ = Petron Live Implementation ( EFG Marketing, San Pablo, Laguna) ~OB Others|Petron Live Implementation (EFG Marketing, San Pablo, Laguna) May 8,9 &amp; 11, 2015| R `""

$ws.Range("P6").Value = "`" ~OT ~ = Home Office Support-San Pablo Laguna ~OB Others|Petron Live Implementation (EFG Marketing, San Pablo, Laguna) May 8,9 &amp; 11, 2015| R `""

$ws.Range("P7").Value = "`" ~OB Others|Petron Live Implementation (EFG Marketing, San Pablo, Laguna) May 8,9 &amp; 11, 2015| R `""

$ws.Range("P14").Value = "`" ~OT ~ = Petron Live Implementation-Bacoor Cavite ~OB Others|Petron Live Implementation ( Bacoor Household Center, Bacoor, Cavite) May 18-29, 2015| R `""

$ws.Range("P15").Value = "`" ~OT ~ = Petron Live Implementation-Bacoor Cavite ~OB Others|Petron Live Implementation ( Bacoor Household Center, Bacoor, Cavite) May 18-29, 2015| R `""

$ws.Range("P16").Value = "`" ~OB Others|Petron Live Implementation ( Bacoor Household Center, Bacoor, Cavite) May 18-29, 2015| R `""

# --- 4. Row 23: ACCUMULATED OT resets from 1.5 to 0 ---
$ws.Range("I23").Value = 0
